$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add header labels for the new columns AI1/AJ1 (mirrors existing "Transacción"/"Comentario" headers)
$ws.Range("AI1").Value = "Transacción"
$ws.Range("AJ1").Value = "Comentario"

# Add the new transaction values in AI2/AJ2
$ws.Range("AI2").Value = "06-3040"
$ws.Range("AJ2").Value = "APROBADO"

# Restore the view to the top-left/A2 selection
$ws.Range("A2").Select()
